$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Val)
    $r = $ws.Range($Address)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "34.070.90"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.788.95"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "227.05"
$ws.Range("E5").Value = "  +1.44%  "
Set-TextValue "D6" "0.546"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue "D8" "32.22"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  +4.07%  "
Set-TextValue "D10" "0.0687"
$ws.Range("E10").Value = "  -2.13%  "
Set-TextValue "D11" "0.0940"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "2.046.62"
$ws.Range("E12").Value = "  +0.47%  "
Set-TextValue "D13" "11.35"
$ws.Range("E13").Value = "  +5.30%  "
$ws.Range("D14").Value = "1.790.87"
$ws.Range("E14").Value = "  -0.04%  "
Set-TextValue "D15" "0.623"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "34.053.63"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("E18").Value = "  +0.86%  "
Set-TextValue "D19" "243.69"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "0.0₃0781"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D21" "10.94"
$ws.Range("E21").Value = "  +2.76%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D22" "1.00"
$ws.Range("E22").Value = "  -0.05%  "
Set-TextValue "D23" "4.10"
$ws.Range("E23").Value = "  +0.44%  "
Set-TextValue "D24" "2.04"
$ws.Range("E24").Value = "  -2.62%  "
Set-TextValue "D25" "161.97"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  +2.53%  "
Set-TextValue "D27" "16.27"
$ws.Range("E27").Value = "  +0.14%  "
Set-TextValue "D28" "0.114"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("E29").Value = "  +0.12%  "
Set-TextValue "D30" "1.23"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").Value = "1.410.35"
$ws.Range("E35").Value = "  +1.41%  "
Set-TextValue "D36" "0.647"
$ws.Range("E36").Value = "  +1.13%  "
Set-TextValue "D37" "0.0190"
$ws.Range("E37").Value = "  +2.76%  "
Set-TextValue "D38" "2.37"
$ws.Range("E38").Value = "  +7.48%  "
$ws.Range("E39").Value = "  -0.50%  "
Set-TextValue "D40" "80.53"
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("E41").Value = "  +0.01%  "
Set-TextValue "D42" "0.919"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("E43").Value = "  -0.04%  "
Set-TextValue "D44" "13.25"
$ws.Range("E44").Value = "  +9.05%  "
$ws.Range("D45").Value = "0.0₆0138"
$ws.Range("E45").Value = "  -6.51%  "
Set-TextValue "D46" "0.0508"
$ws.Range("E46").Value = "  +2.27%  "
Set-TextValue "D47" "6.04"
$ws.Range("E47").Value = "  +3.46%  "
$ws.Range("E48").Value = "  -0.54%  "
Set-TextValue "D49" "106.84"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "1.946.90"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("E51").Value = "  +0.12%  "
